$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 971.8125
$ws.Range("J17").Value = 1003.26666
$ws.Range("L17").Value = 3009.79998
$ws.Range("N17").Value = -3345.79998
$ws.Range("H28").Value = 9141.666999999999
$ws.Range("J28").Value = 8283.333000000001
$ws.Range("L28").Value = 8283.333000000001
$ws.Range("N28").Value = -9253.333000000001
$ws.Range("H92").Value = 1274.5834
$ws.Range("I92").Value = 1004.6
$ws.Range("J92").Value = 2624.5
$ws.Range("K92").Value = 1004.6
$ws.Range("L92").Value = 2624.5
$ws.Range("M92").Value = 243.4
$ws.Range("N92").Value = -5120.5
$ws.Range("H97").Value = 633
$ws.Range("J97").Value = 642.5
$ws.Range("L97").Value = 1927.5
$ws.Range("N97").Value = -2919.5
$ws.Range("H112").Value = 2007.1052
$ws.Range("J112").Value = 2113.8235
$ws.Range("L112").Value = 6341.470499999999
$ws.Range("N112").Value = -8557.470499999999
$ws.Range("H129").Value = 866.4773
$ws.Range("I129").Value = 743.375
$ws.Range("J129").Value = 893.8333
$ws.Range("K129").Value = 2230.125
$ws.Range("L129").Value = 2681.4999
$ws.Range("M129").Value = 2769.875
$ws.Range("N129").Value = -12681.4999
$ws.Range("H138").Value = 2806.5286
$ws.Range("I138").Value = 3366
$ws.Range("J138").Value = 2754.0781
$ws.Range("K138").Value = 10098
$ws.Range("L138").Value = 8262.2343
$ws.Range("M138").Value = -4958
$ws.Range("N138").Value = -18542.2343

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("H32").Value = 7969.48
$ws.Range("I32").Value = 6012.5527
$ws.Range("J32").Value = 14166.417
$ws.Range("K32").Value = 6012.5527
$ws.Range("L32").Value = 14166.417
$ws.Range("M32").Value = -5725.5527
$ws.Range("N32").Value = -14740.417
$ws.Range("H61").Value = 62501430
$ws.Range("I61").Value = 90909900
$ws.Range("J61").Value = 2799.8
$ws.Range("K61").Value = 90909900
$ws.Range("L61").Value = 2799.8
$ws.Range("M61").Value = -90909688
$ws.Range("N61").Value = -3223.8
$ws.Range("H74").Value = 3337.7693
$ws.Range("I74").Value = 1730
$ws.Range("K74").Value = 1730
$ws.Range("M74").Value = -856
$ws.Range("H77").Value = 3337.7693
$ws.Range("I77").Value = 1730
$ws.Range("K77").Value = 8650
$ws.Range("M77").Value = -4282
$ws.Range("H132").Value = 2339.9
$ws.Range("I132").Value = 1784.174
$ws.Range("J132").Value = 4165.857
$ws.Range("K132").Value = 5352.522
$ws.Range("L132").Value = 12497.571
$ws.Range("M132").Value = -2822.522
$ws.Range("N132").Value = -17557.571
$ws.Range("H134").Value = 37700
$ws.Range("J134").Value = 37700
$ws.Range("L134").Value = 37700
$ws.Range("N134").Value = -47840
$ws.Range("H136").Value = 62501430
$ws.Range("I136").Value = 90909900
$ws.Range("J136").Value = 2799.8
$ws.Range("K136").Value = 272729700
$ws.Range("L136").Value = 8399.400000000001
$ws.Range("M136").Value = -272727150
$ws.Range("N136").Value = -13499.4
$ws.Range("M5").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H81").Value = 20467
$ws.Range("J81").Value = 20467
$ws.Range("L81").Value = 20467
$ws.Range("N81").Value = -22589
$ws.Range("H84").Value = 20467
$ws.Range("J84").Value = 20467
$ws.Range("L84").Value = 61401
$ws.Range("N84").Value = -72009
$ws.Range("H134").Value = 5418.5
$ws.Range("I134").Value = 971.26666
$ws.Range("K134").Value = 2913.79998
$ws.Range("M134").Value = -378.7999799999998
$ws.Range("M4").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 184.07692
$ws.Range("I7").Value = 94.09999999999999
$ws.Range("J7").Value = 484
$ws.Range("K7").Value = 94.09999999999999
$ws.Range("L7").Value = 484
$ws.Range("M7").Value = 18.90000000000001
$ws.Range("N7").Value = -710
$ws.Range("H58").Value = 6563.85
$ws.Range("I58").Value = 938
$ws.Range("K58").Value = 938
$ws.Range("M58").Value = -735
$ws.Range("H62").Value = 15386840
$ws.Range("I62").Value = 2374.4546
$ws.Range("K62").Value = 2374.4546
$ws.Range("M62").Value = -1750.4546
$ws.Range("H65").Value = 15386840
$ws.Range("I65").Value = 2374.4546
$ws.Range("K65").Value = 11872.273
$ws.Range("M65").Value = -8752.273000000001
$ws.Range("H81").Value = 16500
$ws.Range("J81").Value = 16500
$ws.Range("L81").Value = 16500
$ws.Range("N81").Value = -18496
$ws.Range("H84").Value = 16500
$ws.Range("J84").Value = 16500
$ws.Range("L84").Value = 49500
$ws.Range("N84").Value = -59484
$ws.Range("H136").Value = 6563.85
$ws.Range("I136").Value = 938
$ws.Range("K136").Value = 2814
$ws.Range("M136").Value = -264
$ws.Range("H141").Value = 296790.72
$ws.Range("J141").Value = 296790.72
$ws.Range("L141").Value = 296790.72
$ws.Range("N141").Value = -307150.72

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1263
$ws.Range("J87").Value = 2001
$ws.Range("L87").Value = 6003
$ws.Range("N87").Value = -8499
$ws.Range("H90").Value = 1263
$ws.Range("J90").Value = 2001
$ws.Range("L90").Value = 18009
$ws.Range("N90").Value = -30489
$ws.Range("H131").Value = 27068552
$ws.Range("J131").Value = 54728.715
$ws.Range("L131").Value = 164186.145
$ws.Range("N131").Value = -174266.145
$ws.Range("H140").Value = 28976.447
$ws.Range("I140").Value = 57992.723
$ws.Range("J140").Value = 2861.8
$ws.Range("K140").Value = 173978.169
$ws.Range("L140").Value = 8585.400000000001
$ws.Range("M140").Value = -168798.169
$ws.Range("N140").Value = -18945.4
$ws.Range("H141").Value = 55558960
$ws.Range("I141").Value = 62502740
$ws.Range("K141").Value = 187508220
$ws.Range("M141").Value = -187503040

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 4000
$ws.Range("J47").Value = 4000
$ws.Range("L47").Value = 4000
$ws.Range("N47").Value = -5136
$ws.Range("H123").Value = 22488
$ws.Range("J123").Value = 22488
$ws.Range("L123").Value = 22488
$ws.Range("N123").Value = -27388
$ws.Range("H132").Value = 5074.9473
$ws.Range("I132").Value = 5997.25
$ws.Range("J132").Value = 3493.8572
$ws.Range("K132").Value = 17991.75
$ws.Range("L132").Value = 10481.5716
$ws.Range("M132").Value = -15461.75
$ws.Range("N132").Value = -15541.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1988.5
$ws.Range("I7").Value = 1625.7142
$ws.Range("K7").Value = 1625.7142
$ws.Range("M7").Value = -1513.7142
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("H126").Value = 1988.5
$ws.Range("I126").Value = 1625.7142
$ws.Range("K126").Value = 4877.142599999999
$ws.Range("M126").Value = -2407.142599999999
$ws.Range("H136").Value = 1796.1666
$ws.Range("I136").Value = 1743.826
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5231.478
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2681.478
$ws.Range("N136").Value = -14100
$ws.Range("M26").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 64996.668
$ws.Range("J135").Value = 64996.668
$ws.Range("L135").Value = 64996.668
$ws.Range("N135").Value = -75136.66800000001
$ws.Range("H136").Value = 1582.6444
$ws.Range("I136").Value = 698.5238000000001
$ws.Range("J136").Value = 2356.25
$ws.Range("K136").Value = 2095.5714
$ws.Range("L136").Value = 7068.75
$ws.Range("M136").Value = 454.4285999999997
$ws.Range("N136").Value = -12168.75
